$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @("codeforiati:category-name","codeforiati:group-name","codeforiati:category-code","codeforiati:group-code"),
    @("Education, Level Unspecified","Education","111","110"),
    @("Education, Level Unspecified","Education","111","110"),
    @("Education, Level Unspecified","Education","111","110"),
    @("Education, Level Unspecified","Education","111","110"),
    @("Basic Education","Education","112","110"),
    @("Basic Education","Education","112","110"),
    @("Basic Education","Education","112","110"),
    @("Basic Education","Education","112","110"),
    @("Basic Education","Education","112","110"),
    @("Basic Education","Education","112","110"),
    @("Basic Education","Education","112","110"),
    @("Secondary Education","Education","113","110"),
    @("Secondary Education","Education","113","110"),
    @("Post-Secondary Education","Education","114","110"),
    @("Post-Secondary Education","Education","114","110"),
    @("Health, General","Health","121","120"),
    @("Health, General","Health","121","120"),
    @("Health, General","Health","121","120"),
    @("Health, General","Health","121","120"),
    @("Basic Health","Health","122","120"),
    @("Basic Health","Health","122","120"),
    @("Basic Health","Health","122","120"),
    @("Basic Health","Health","122","120"),
    @("Basic Health","Health","122","120"),
    @("Basic Health","Health","122","120"),
    @("Basic Health","Health","122","120"),
    @("Basic Health","Health","122","120"),
    @("Basic Health","Health","122","120"),
    @("Non-communicable diseases (NCDs)","Health","123","120"),
    @("Non-communicable diseases (NCDs)","Health","123","120"),
    @("Non-communicable diseases (NCDs)","Health","123","120"),
    @("Non-communicable diseases (NCDs)","Health","123","120"),
    @("Non-communicable diseases (NCDs)","Health","123","120"),
    @("Non-communicable diseases (NCDs)","Health","123","120"),
    @("Population Policies/Programmes & Reproductive Health","Population Policies/Programmes & Reproductive Health","130","130"),
    @("Population Policies/Programmes & Reproductive Health","Population Policies/Programmes & Reproductive Health","130","130"),
    @("Population Policies/Programmes & Reproductive Health","Population Policies/Programmes & Reproductive Health","130","130"),
    @("Population Policies/Programmes & Reproductive Health","Population Policies/Programmes & Reproductive Health","130","130"),
    @("Population Policies/Programmes & Reproductive Health","Population Policies/Programmes & Reproductive Health","130","130"),
    @("Water Supply & Sanitation","Water Supply & Sanitation","140","140"),
    @("Water Supply & Sanitation","Water Supply & Sanitation","140","140"),
    @("Water Supply & Sanitation","Water Supply & Sanitation","140","140"),
    @("Water Supply & Sanitation","Water Supply & Sanitation","140","140"),
    @("Water Supply & Sanitation","Water Supply & Sanitation","140","140"),
    @("Water Supply & Sanitation","Water Supply & Sanitation","140","140"),
    @("Water Supply & Sanitation","Water Supply & Sanitation","140","140"),
    @("Water Supply & Sanitation","Water Supply & Sanitation","140","140"),
    @("Water Supply & Sanitation","Water Supply & Sanitation","140","140"),
    @("Water Supply & Sanitation","Water Supply & Sanitation","140","140"),
    @("Water Supply & Sanitation","Water Supply & Sanitation","140","140"),
    @("Government & Civil Society-general","Government & Civil Society","151","150"),
    @("Government & Civil Society-general","Government & Civil Society","151","150"),
    @("Government & Civil Society-general","Government & Civil Society","151","150"),
    @("Government & Civil Society-general","Government & Civil Society","151","150"),
    @("Government & Civil Society-general","Government & Civil Society","151","150"),
    @("Government & Civil Society-general","Government & Civil Society","151","150"),
    @("Government & Civil Society-general","Government & Civil Society","151","150"),
    @("Government & Civil Society-general","Government & Civil Society","151","150"),
    @("Government & Civil Society-general","Government & Civil Society","151","150"),
    @("Government & Civil Society-general","Government & Civil Society","151","150"),
    @("Government & Civil Society-general","Government & Civil Society","151","150"),
    @("Government & Civil Society-general","Government & Civil Society","151","150"),
    @("Government & Civil Society-general","Government & Civil Society","151","150"),
    @("Government & Civil Society-general","Government & Civil Society","151","150"),
    @("Government & Civil Society-general","Government & Civil Society","151","150"),
    @("Government & Civil Society-general","Government & Civil Society","151","150"),
    @("Conflict, Peace & Security","Government & Civil Society","152","150"),
    @("Conflict, Peace & Security","Government & Civil Society","152","150"),
    @("Conflict, Peace & Security","Government & Civil Society","152","150"),
    @("Conflict, Peace & Security","Government & Civil Society","152","150"),
    @("Conflict, Peace & Security","Government & Civil Society","152","150"),
    @("Conflict, Peace & Security","Government & Civil Society","152","150"),
    @("Other Social Infrastructure & Services","Other Social Infrastructure & Services","160","160"),
    @("Other Social Infrastructure & Services","Other Social Infrastructure & Services","160","160"),
    @("Other Social Infrastructure & Services","Other Social Infrastructure & Services","160","160"),
    @("Other Social Infrastructure & Services","Other Social Infrastructure & Services","160","160"),
    @("Other Social Infrastructure & Services","Other Social Infrastructure & Services","160","160"),
    @("Other Social Infrastructure & Services","Other Social Infrastructure & Services","160","160"),
    @("Other Social Infrastructure & Services","Other Social Infrastructure & Services","160","160"),
    @("Other Social Infrastructure & Services","Other Social Infrastructure & Services","160","160"),
    @("Other Social Infrastructure & Services","Other Social Infrastructure & Services","160","160"),
    @("Other Social Infrastructure & Services","Other Social Infrastructure & Services","160","160"),
    @("Other Social Infrastructure & Services","Other Social Infrastructure & Services","160","160"),
    @("Transport & Storage","Transport & Storage","210","210"),
    @("Transport & Storage","Transport & Storage","210","210"),
    @("Transport & Storage","Transport & Storage","210","210"),
    @("Transport & Storage","Transport & Storage","210","210"),
    @("Transport & Storage","Transport & Storage","210","210"),
    @("Transport & Storage","Transport & Storage","210","210"),
    @("Transport & Storage","Transport & Storage","210","210"),
    @("Communications","Communications","220","220"),
    @("Communications","Communications","220","220"),
    @("Communications","Communications","220","220"),
    @("Communications","Communications","220","220"),
    @("Energy Policy","Energy","231","230"),
    @("Energy Policy","Energy","231","230"),
    @("Energy Policy","Energy","231","230"),
    @("Energy Policy","Energy","231","230"),
    @("Energy generation, renewable sources","Energy","232","230"),
    @("Energy generation, renewable sources","Energy","232","230"),
    @("Energy generation, renewable sources","Energy","232","230"),
    @("Energy generation, renewable sources","Energy","232","230"),
    @("Energy generation, renewable sources","Energy","232","230"),
    @("Energy generation, renewable sources","Energy","232","230"),
    @("Energy generation, renewable sources","Energy","232","230"),
    @("Energy generation, renewable sources","Energy","232","230"),
    @("Energy generation, renewable sources","Energy","232","230"),
    @("Energy generation, non-renewable sources","Energy","233","230"),
    @("Energy generation, non-renewable sources","Energy","233","230"),
    @("Energy generation, non-renewable sources","Energy","233","230"),
    @("Energy generation, non-renewable sources","Energy","233","230"),
    @("Energy generation, non-renewable sources","Energy","233","230"),
    @("Energy generation, non-renewable sources","Energy","233","230"),
    @("Hybrid energy plants","Energy","234","230"),
    @("Nuclear energy plants","Energy","235","230"),
    @("Energy distribution","Energy","236","230"),
    @("Energy distribution","Energy","236","230"),
    @("Energy distribution","Energy","236","230"),
    @("Energy distribution","Energy","236","230"),
    @("Energy distribution","Energy","236","230"),
    @("Energy distribution","Energy","236","230"),
    @("Energy distribution","Energy","236","230"),
    @("Banking & Financial Services","Banking & Financial Services","240","240"),
    @("Banking & Financial Services","Banking & Financial Services","240","240"),
    @("Banking & Financial Services","Banking & Financial Services","240","240"),
    @("Banking & Financial Services","Banking & Financial Services","240","240"),
    @("Banking & Financial Services","Banking & Financial Services","240","240"),
    @("Banking & Financial Services","Banking & Financial Services","240","240"),
    @("Business & Other Services","Business & Other Services","250","250"),
    @("Business & Other Services","Business & Other Services","250","250"),
    @("Business & Other Services","Business & Other Services","250","250"),
    @("Business & Other Services","Business & Other Services","250","250"),
    @("Agriculture","Agriculture, Forestry, Fishing","311","310"),
    @("Agriculture","Agriculture, Forestry, Fishing","311","310"),
    @("Agriculture","Agriculture, Forestry, Fishing","311","310"),
    @("Agriculture","Agriculture, Forestry, Fishing","311","310"),
    @("Agriculture","Agriculture, Forestry, Fishing","311","310"),
    @("Agriculture","Agriculture, Forestry, Fishing","311","310"),
    @("Agriculture","Agriculture, Forestry, Fishing","311","310"),
    @("Agriculture","Agriculture, Forestry, Fishing","311","310"),
    @("Agriculture","Agriculture, Forestry, Fishing","311","310"),
    @("Agriculture","Agriculture, Forestry, Fishing","311","310"),
    @("Agriculture","Agriculture, Forestry, Fishing","311","310"),
    @("Agriculture","Agriculture, Forestry, Fishing","311","310"),
    @("Agriculture","Agriculture, Forestry, Fishing","311","310"),
    @("Agriculture","Agriculture, Forestry, Fishing","311","310"),
    @("Agriculture","Agriculture, Forestry, Fishing","311","310"),
    @("Agriculture","Agriculture, Forestry, Fishing","311","310"),
    @("Agriculture","Agriculture, Forestry, Fishing","311","310"),
    @("Agriculture","Agriculture, Forestry, Fishing","311","310"),
    @("Forestry","Agriculture, Forestry, Fishing","312","310"),
    @("Forestry","Agriculture, Forestry, Fishing","312","310"),
    @("Forestry","Agriculture, Forestry, Fishing","312","310"),
    @("Forestry","Agriculture, Forestry, Fishing","312","310"),
    @("Forestry","Agriculture, Forestry, Fishing","312","310"),
    @("Forestry","Agriculture, Forestry, Fishing","312","310"),
    @("Fishing","Agriculture, Forestry, Fishing","313","310"),
    @("Fishing","Agriculture, Forestry, Fishing","313","310"),
    @("Fishing","Agriculture, Forestry, Fishing","313","310"),
    @("Fishing","Agriculture, Forestry, Fishing","313","310"),
    @("Fishing","Agriculture, Forestry, Fishing","313","310"),
    @("Industry","Industry, Mining, Construction","321","320"),
    @("Industry","Industry, Mining, Construction","321","320"),
    @("Industry","Industry, Mining, Construction","321","320"),
    @("Industry","Industry, Mining, Construction","321","320"),
    @("Industry","Industry, Mining, Construction","321","320"),
    @("Industry","Industry, Mining, Construction","321","320"),
    @("Industry","Industry, Mining, Construction","321","320"),
    @("Industry","Industry, Mining, Construction","321","320"),
    @("Industry","Industry, Mining, Construction","321","320"),
    @("Industry","Industry, Mining, Construction","321","320"),
    @("Industry","Industry, Mining, Construction","321","320"),
    @("Industry","Industry, Mining, Construction","321","320"),
    @("Industry","Industry, Mining, Construction","321","320"),
    @("Industry","Industry, Mining, Construction","321","320"),
    @("Industry","Industry, Mining, Construction","321","320"),
    @("Industry","Industry, Mining, Construction","321","320"),
    @("Industry","Industry, Mining, Construction","321","320"),
    @("Industry","Industry, Mining, Construction","321","320"),
    @("Industry","Industry, Mining, Construction","321","320"),
    @("Mineral Resources & Mining","Industry, Mining, Construction","322","320"),
    @("Mineral Resources & Mining","Industry, Mining, Construction","322","320"),
    @("Mineral Resources & Mining","Industry, Mining, Construction","322","320"),
    @("Mineral Resources & Mining","Industry, Mining, Construction","322","320"),
    @("Mineral Resources & Mining","Industry, Mining, Construction","322","320"),
    @("Mineral Resources & Mining","Industry, Mining, Construction","322","320"),
    @("Mineral Resources & Mining","Industry, Mining, Construction","322","320"),
    @("Mineral Resources & Mining","Industry, Mining, Construction","322","320"),
    @("Mineral Resources & Mining","Industry, Mining, Construction","322","320"),
    @("Mineral Resources & Mining","Industry, Mining, Construction","322","320"),
    @("Construction","Industry, Mining, Construction","323","320"),
    @("Trade Policies & Regulations","Trade Policies & Regulations","331","331"),
    @("Trade Policies & Regulations","Trade Policies & Regulations","331","331"),
    @("Trade Policies & Regulations","Trade Policies & Regulations","331","331"),
    @("Trade Policies & Regulations","Trade Policies & Regulations","331","331"),
    @("Trade Policies & Regulations","Trade Policies & Regulations","331","331"),
    @("Trade Policies & Regulations","Trade Policies & Regulations","331","331"),
    @("Tourism","Tourism","332","332"),
    @("General Environment Protection","General Environment Protection","410","410"),
    @("General Environment Protection","General Environment Protection","410","410"),
    @("General Environment Protection","General Environment Protection","410","410"),
    @("General Environment Protection","General Environment Protection","410","410"),
    @("General Environment Protection","General Environment Protection","410","410"),
    @("General Environment Protection","General Environment Protection","410","410"),
    @("Other Multisector","Other Multisector","430","430"),
    @("Other Multisector","Other Multisector","430","430"),
    @("Other Multisector","Other Multisector","430","430"),
    @("Other Multisector","Other Multisector","430","430"),
    @("Other Multisector","Other Multisector","430","430"),
    @("Other Multisector","Other Multisector","430","430"),
    @("Other Multisector","Other Multisector","430","430"),
    @("Other Multisector","Other Multisector","430","430"),
    @("Other Multisector","Other Multisector","430","430"),
    @("Other Multisector","Other Multisector","430","430"),
    @("General Budget Support","General Budget Support","510","510"),
    @("Development Food Assistance","Development Food Assistance","520","520"),
    @("Other Commodity Assistance","Other Commodity Assistance","530","530"),
    @("Other Commodity Assistance","Other Commodity Assistance","530","530"),
    @("Action Relating to Debt","Action Relating to Debt","600","600"),
    @("Action Relating to Debt","Action Relating to Debt","600","600"),
    @("Action Relating to Debt","Action Relating to Debt","600","600"),
    @("Action Relating to Debt","Action Relating to Debt","600","600"),
    @("Action Relating to Debt","Action Relating to Debt","600","600"),
    @("Action Relating to Debt","Action Relating to Debt","600","600"),
    @("Action Relating to Debt","Action Relating to Debt","600","600"),
    @("Emergency Response","Emergency Response","720","720"),
    @("Emergency Response","Emergency Response","720","720"),
    @("Emergency Response","Emergency Response","720","720"),
    @("Reconstruction Relief & Rehabilitation","Reconstruction Relief & Rehabilitation","730","730"),
    @("Disaster Prevention & Preparedness","Disaster Prevention & Preparedness","740","740"),
    @("Administrative Costs of Donors","Administrative Costs of Donors","910","910"),
    @("Refugees in Donor Countries","Refugees in Donor Countries","930","930"),
    @("Unallocated / Unspecified","Unallocated / Unspecified","998","998"),
    @("Unallocated / Unspecified","Unallocated / Unspecified","998","998")
)

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $i + 1
    $vals = $rows[$i]
    $ws.Cells.Item($r, 4).Value2 = $vals[0]
    $ws.Cells.Item($r, 5).Value2 = $vals[1]
    $ws.Cells.Item($r, 6).Value2 = $vals[2]
    $ws.Cells.Item($r, 7).Value2 = $vals[3]
}

Write-Host "Done updating" $rows.Count "rows"
